# Apply "repull data, push all data, mean calculation" update.
# Column F ("dSF") values were re-pulled / recalculated for several rows;
# only the F column differs from the previous snapshot, per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -2
    10 = -3
    11 = 9
    12 = 8
    13 = -10
    15 = -3
    17 = 0
    24 = 4
    27 = -2
    29 = 0
    35 = 2
    37 = -1
    39 = -1
    40 = 0
    44 = 3
    46 = 1
    49 = -1
    53 = -1
    54 = 1
    55 = 2
    56 = -4
    59 = -1
    61 = -1
    69 = -1
    71 = 4
    72 = 5
    74 = 0
    81 = 3
    92 = -7
    93 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

Write-Host "Updated $($updates.Count) cells in column F"
